$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New working set of 192 German verb sequences (column B, rows 2-193)
$words = New-Object 'object[,]' 192,1
$words[0,0] = 'ehren'
$words[1,0] = 'fließen'
$words[2,0] = 'ändern'
$words[3,0] = 'sichern'
$words[4,0] = 'trennen'
$words[5,0] = 'steuern'
$words[6,0] = 'filmen'
$words[7,0] = 'malen'
$words[8,0] = 'rufen'
$words[9,0] = 'schwingen'
$words[10,0] = 'teilen'
$words[11,0] = 'formen'
$words[12,0] = 'flehen'
$words[13,0] = 'ziehen'
$words[14,0] = 'pflegen'
$words[15,0] = 'helfen'
$words[16,0] = 'jagen'
$words[17,0] = 'fallen'
$words[18,0] = 'heilen'
$words[19,0] = 'ächzen'
$words[20,0] = 'pfeifen'
$words[21,0] = 'nähen'
$words[22,0] = 'lockern'
$words[23,0] = 'streichen'
$words[24,0] = 'foltern'
$words[25,0] = 'geben'
$words[26,0] = 'spielen'
$words[27,0] = 'greifen'
$words[28,0] = 'wehtun'
$words[29,0] = 'herrschen'
$words[30,0] = 'treten'
$words[31,0] = 'töten'
$words[32,0] = 'ärgern'
$words[33,0] = 'zünden'
$words[34,0] = 'dienen'
$words[35,0] = 'drehen'
$words[36,0] = 'sparen'
$words[37,0] = 'graben'
$words[38,0] = 'regnen'
$words[39,0] = 'wundern'
$words[40,0] = 'feiern'
$words[41,0] = 'tropfen'
$words[42,0] = 'irren'
$words[43,0] = 'kürzen'
$words[44,0] = 'garen'
$words[45,0] = 'sterben'
$words[46,0] = 'segnen'
$words[47,0] = 'freuen'
$words[48,0] = 'klingen'
$words[49,0] = 'läuten'
$words[50,0] = 'schwächen'
$words[51,0] = 'werden'
$words[52,0] = 'orten'
$words[53,0] = 'siegen'
$words[54,0] = 'fällen'
$words[55,0] = 'flüchten'
$words[56,0] = 'opfern'
$words[57,0] = 'loben'
$words[58,0] = 'treiben'
$words[59,0] = 'dringen'
$words[60,0] = 'hören'
$words[61,0] = 'spinnen'
$words[62,0] = 'bluten'
$words[63,0] = 'hupen'
$words[64,0] = 'bilden'
$words[65,0] = 'wahren'
$words[66,0] = 'backen'
$words[67,0] = 'testen'
$words[68,0] = 'enden'
$words[69,0] = 'planen'
$words[70,0] = 'spenden'
$words[71,0] = 'boxen'
$words[72,0] = 'liefern'
$words[73,0] = 'platzen'
$words[74,0] = 'kaufen'
$words[75,0] = 'betteln'
$words[76,0] = 'fischen'
$words[77,0] = 'fügen'
$words[78,0] = 'achten'
$words[79,0] = 'bitten'
$words[80,0] = 'küssen'
$words[81,0] = 'schwören'
$words[82,0] = 'biegen'
$words[83,0] = 'heulen'
$words[84,0] = 'klettern'
$words[85,0] = 'mauern'
$words[86,0] = 'starren'
$words[87,0] = 'schlucken'
$words[88,0] = 'spüren'
$words[89,0] = 'warnen'
$words[90,0] = 'knien'
$words[91,0] = 'münzen'
$words[92,0] = 'reisen'
$words[93,0] = 'quälen'
$words[94,0] = 'klagen'
$words[95,0] = 'pflanzen'
$words[96,0] = 'posten'
$words[97,0] = 'schultern'
$words[98,0] = 'räumen'
$words[99,0] = 'suchen'
$words[100,0] = 'sprengen'
$words[101,0] = 'hassen'
$words[102,0] = 'wüten'
$words[103,0] = 'stehlen'
$words[104,0] = 'zeigen'
$words[105,0] = 'schalten'
$words[106,0] = 'ahnen'
$words[107,0] = 'kehren'
$words[108,0] = 'regeln'
$words[109,0] = 'baden'
$words[110,0] = 'wachsen'
$words[111,0] = 'sorgen'
$words[112,0] = 'lügen'
$words[113,0] = 'parken'
$words[114,0] = 'runden'
$words[115,0] = 'zögern'
$words[116,0] = 'gelten'
$words[117,0] = 'sprechen'
$words[118,0] = 'beißen'
$words[119,0] = 'lesen'
$words[120,0] = 'hauen'
$words[121,0] = 'messen'
$words[122,0] = 'tollen'
$words[123,0] = 'machen'
$words[124,0] = 'scheitern'
$words[125,0] = 'knarren'
$words[126,0] = 'schleppen'
$words[127,0] = 'sinken'
$words[128,0] = 'bergen'
$words[129,0] = 'duschen'
$words[130,0] = 'heben'
$words[131,0] = 'schulden'
$words[132,0] = 'kümmern'
$words[133,0] = 'leugnen'
$words[134,0] = 'träumen'
$words[135,0] = 'jubeln'
$words[136,0] = 'stimmen'
$words[137,0] = 'nutzen'
$words[138,0] = 'sperren'
$words[139,0] = 'wagen'
$words[140,0] = 'grüßen'
$words[141,0] = 'saufen'
$words[142,0] = 'trauen'
$words[143,0] = 'segeln'
$words[144,0] = 'arten'
$words[145,0] = 'rasen'
$words[146,0] = 'reizen'
$words[147,0] = 'fahren'
$words[148,0] = 'schrecken'
$words[149,0] = 'äußern'
$words[150,0] = 'lösen'
$words[151,0] = 'bellen'
$words[152,0] = 'wirken'
$words[153,0] = 'hindern'
$words[154,0] = 'schenken'
$words[155,0] = 'scheinen'
$words[156,0] = 'liegen'
$words[157,0] = 'wechseln'
$words[158,0] = 'strahlen'
$words[159,0] = 'buchen'
$words[160,0] = 'stammen'
$words[161,0] = 'bleiben'
$words[162,0] = 'kosten'
$words[163,0] = 'folgen'
$words[164,0] = 'wenden'
$words[165,0] = 'lohnen'
$words[166,0] = 'altern'
$words[167,0] = 'zielen'
$words[168,0] = 'brauchen'
$words[169,0] = 'drücken'
$words[170,0] = 'bauen'
$words[171,0] = 'schreiten'
$words[172,0] = 'tragen'
$words[173,0] = 'stechen'
$words[174,0] = 'trotzen'
$words[175,0] = 'schmecken'
$words[176,0] = 'leeren'
$words[177,0] = 'wärmen'
$words[178,0] = 'gründen'
$words[179,0] = 'stecken'
$words[180,0] = 'bieten'
$words[181,0] = 'weichen'
$words[182,0] = 'kichern'
$words[183,0] = 'werfen'
$words[184,0] = 'fangen'
$words[185,0] = 'schließen'
$words[186,0] = 'decken'
$words[187,0] = 'führen'
$words[188,0] = 'seufzen'
$words[189,0] = 'erben'
$words[190,0] = 'mögen'
$words[191,0] = 'wohnen'

$ws.Range("B2:B193").Value = $words

